# Apply updated crypto price/volume data (GitHub Actions refresh).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "70.312.50"
$ws.Range("D3").Value = "3.510.37"
$ws.Range("E3").Value = "  +0.29%  "
$ws.Range("E4").Value = "  -0.08%  "
$ws.Range("D5").Value = "'604.98"
$ws.Range("E5").Value = "  -0.12%  "
$ws.Range("D6").Value = "'175.00"
$ws.Range("E6").Value = "  +3.15%  "
$ws.Range("E7").Value = "  -1.00%  "
$ws.Range("D8").Value = "3.505.28"
$ws.Range("E8").Value = "  +0.26%  "
$ws.Range("E9").Value = "  +0.00%  "
$ws.Range("E10").Value = "  -1.96%  "
$ws.Range("D11").Value = "'7.17"
$ws.Range("E11").Value = "  +7.58%  "
$ws.Range("D12").Value = "'0.583"
$ws.Range("E12").Value = "  +0.77%  "
$ws.Range("D13").Value = "'46.19"
$ws.Range("E13").Value = "  -1.59%  "
$ws.Range("E14").Value = "  -1.01%  "
$ws.Range("D15").Value = "4.070.08"
$ws.Range("E15").Value = "  -0.27%  "
$ws.Range("B16").Value = "BitcoinCash"
$ws.Range("C16").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Range("D16").Value = "'612.81"
$ws.Range("E16").Value = "  +0.26%  "
$ws.Range("B17").Value = "Polkadot"
$ws.Range("C17").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D17").Value = "'8.28"
$ws.Range("E17").Value = "  -0.06%  "
$ws.Range("D18").Value = "3.509.92"
$ws.Range("E18").Value = "  +0.22%  "
$ws.Range("D19").Value = "70.441.75"
$ws.Range("E19").Value = "  +1.31%  "
$ws.Range("E20").Value = "  +1.04%  "
$ws.Range("D21").Value = "'17.35"
$ws.Range("E21").Value = "  +0.88%  "
$ws.Range("E22").Value = "  -0.39%  "
$ws.Range("D23").Value = "'9.01"
$ws.Range("E23").Value = "  -11.34%  "
$ws.Range("D24").Value = "'98.48"
$ws.Range("E24").Value = "  +3.14%  "
$ws.Range("E26").Value = "  -3.49%  "
$ws.Range("E27").Value = "  -0.12%  "
$ws.Range("E28").Value = "  -1.29%  "
$ws.Range("D29").Value = "'33.86"
$ws.Range("E29").Value = "  +2.40%  "
$ws.Range("E30").Value = "  -1.85%  "
$ws.Range("D31").Value = "'2.97"
$ws.Range("E31").Value = "  -3.32%  "
$ws.Range("D32").Value = "'8.02"
$ws.Range("E32").Value = "  -4.50%  "
$ws.Range("B33").Value = "Bittensor"
$ws.Range("C33").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D33").Value = "'638.38"
$ws.Range("E33").Value = "  +15.29%  "
$ws.Range("B34").Value = "Mantle"
$ws.Range("C34").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D34").Value = "'1.28"
$ws.Range("E34").Value = "  -3.58%  "
$ws.Range("B35").Value = "NEARProtocol"
$ws.Range("C35").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D35").Value = "'6.82"
$ws.Range("E35").Value = "  -1.02%  "
$ws.Range("E36").Value = "  -1.35%  "
$ws.Range("E37").Value = "  +2.17%  "
$ws.Range("E38").Value = "  +0.25%  "
$ws.Range("D39").Value = "'0.0474"
$ws.Range("E39").Value = "  +5.82%  "
$ws.Range("D40").Value = "'56.77"
$ws.Range("E41").Value = "  +0.05%  "
$ws.Range("E42").Value = "  +2.27%  "
$ws.Range("D43").Value = "3.367.74"
$ws.Range("E43").Value = "  +1.16%  "
$ws.Range("E44").Value = "  +6.06%  "
$ws.Range("E45").Value = "  -5.08%  "
$ws.Range("D46").Value = "'32.14"
$ws.Range("E46").Value = "  -2.35%  "
$ws.Range("D47").Value = "'2.89"
$ws.Range("E47").Value = "  -0.15%  "
$ws.Range("D48").Value = "'2.55"
$ws.Range("E48").Value = "  -1.70%  "
$ws.Range("E49").Value = "  +0.78%  "
$ws.Range("D50").Value = "'132.60"
$ws.Range("E50").Value = "  -2.38%  "
